$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 1-9, 11-12 are untouched by this edit; only row 10 (B/C text) and
# rows 13-23 (reshuffled "Programa"/"Docentes"/"Bibliografia" block) change.

# --- Row 10: only the Objetivos (B/C) text changes; keep style & height. ---
$ws.Range("B10").Value = 'Dotar os alunos dos conhecimentos de biologia celular abrangendo a organização estrutural e molecular da célula, proporcionando os conhecimentos básicos de biologia geral necessários à compreensão das demais disciplinas correlacionadas do curso de Engenharia Bioquímica.'
$ws.Range("C10").Value = 'Dotar os alunos dos conhecimentos de biologia celular abrangendo a organização estrutural e molecular da célula, proporcionando os conhecimentos básicos de biologia geral necessários à compreensão das demais disciplinas correlacionadas do curso de Engenharia Bioquímica.'

# --- Save format templates (col A/B/C canonical styles) into a scratch row ---
$ws.Range("A3").Copy($ws.Range("A200"))
$ws.Range("B3").Copy($ws.Range("B200"))
$ws.Range("C3").Copy($ws.Range("C200"))

# --- Wipe rows 13-23 (content + formatting) and reset row heights to default ---
$ws.Range("A13:C23").Clear()
for ($r = 13; $r -le 23; $r++) {
    $ws.Rows.Item($r).AutoFit()
}

# --- Rebuild rows 13-23 with the target content ---
# Row 13
$ws.Range("B200").Copy($ws.Range("B13"))
$ws.Range("B13").Value = '1304060 - Maria das Graças de Almeida Felipe'
$ws.Range("C200").Copy($ws.Range("C13"))
$ws.Range("C13").Value = '1304060 - Maria das Graças de Almeida Felipe'

# Row 14
$ws.Range("B200").Copy($ws.Range("B14"))
$ws.Range("B14").Value = '8853480 - Tatiane da Franca Silva'
$ws.Range("C200").Copy($ws.Range("C14"))
$ws.Range("C14").Value = '8853480 - Tatiane da Franca Silva'

# Row 15
$ws.Range("A200").Copy($ws.Range("A15"))
$ws.Range("A15").Value = 'Programa resumido:'
$ws.Range("B200").Copy($ws.Range("B15"))
$ws.Range("B15").Value = 'Origem e evolução das células; análise estrutural das células; organização interna das células.'
$ws.Range("C200").Copy($ws.Range("C15"))
$ws.Range("C15").Value = 'Origem e evolução das células; análise estrutural das células; organização interna das células.'
$ws.Rows.Item(15).RowHeight = 60

# Row 16
$ws.Range("A200").Copy($ws.Range("A16"))
$ws.Range("A16").Value = 'Short syllabus:'
$ws.Range("B200").Copy($ws.Range("B16"))
$ws.Range("B16").Value = 'Cellular origin and evolution; structural analysis of cells; internal organization of cells.'
$ws.Range("C200").Copy($ws.Range("C16"))
$ws.Range("C16").Value = 'Cellular origin and evolution; structural analysis of cells; internal organization of cells.'
$ws.Rows.Item(16).RowHeight = 60

# Row 17
$ws.Range("A200").Copy($ws.Range("A17"))
$ws.Range("A17").Value = 'Programa:'
$ws.Range("B200").Copy($ws.Range("B17"))
$ws.Range("B17").Value = 'Origem e evolução das células: Conceitos básicos de sistemática e filogenia molecular, características dos três domínios. –Análise estrutural das células ao microscópio: Microscopia ótica e microscopia eletrônica.–Organização interna das células: Células procarióticas e eucarióticas; estrutura e transporte através das membranas; compartimentos intracelulares (núcleo, retículo endoplasmático rugoso e liso, complexo de golgi, lisossomos e peroxissomos) e endereçamento de proteínas; tráfego intracelular de vesículas (via secretora e endocítica); conversão de energia (mitocôndria e cloroplasto); comunicação e sinalização celular; citoesqueleto; ciclo e divisão celular (mitose e meiose); matriz extracelulares e parede celular vegetal.'
$ws.Range("C200").Copy($ws.Range("C17"))
$ws.Range("C17").Value = 'Origem e evolução das células: Conceitos básicos de sistemática e filogenia molecular, características dos três domínios. –Análise estrutural das células ao microscópio: Microscopia ótica e microscopia eletrônica.–Organização interna das células: Células procarióticas e eucarióticas; estrutura e transporte através das membranas; compartimentos intracelulares (núcleo, retículo endoplasmático rugoso e liso, complexo de golgi, lisossomos e peroxissomos) e endereçamento de proteínas; tráfego intracelular de vesículas (via secretora e endocítica); conversão de energia (mitocôndria e cloroplasto); comunicação e sinalização celular; citoesqueleto; ciclo e divisão celular (mitose e meiose); matriz extracelulares e parede celular vegetal.'
$ws.Rows.Item(17).RowHeight = 120

# Row 18
$ws.Range("A200").Copy($ws.Range("A18"))
$ws.Range("A18").Value = 'Syllabus:'
$ws.Range("B200").Copy($ws.Range("B18"))
$ws.Range("B18").Value = 'Origin and evolution of cells: basic concepts of systematic and molecular phylogeny, characteristics of the three domains. Structural analysis of cells at the microscope: optical and electronic microscopy.Internal organization of cells:Prokaryotic and eukaryotic cells; structure and transport through the membranes; intracelular compartments (nucleus, rough and smooth endoplasmic reticulum, golgi complex; lysosomes and peroxisomes) and protein addressing; intracellular traffic of vesicles (secretory and endocytic pathway); energy conversion (mitochondria and chloroplast); cellular communication and signalization; cytoskeleton; cellular cycle and division (mitosis and meiosis); extracellular matrix and vegetal cell wall.'
$ws.Range("C200").Copy($ws.Range("C18"))
$ws.Range("C18").Value = 'Origin and evolution of cells: basic concepts of systematic and molecular phylogeny, characteristics of the three domains. Structural analysis of cells at the microscope: optical and electronic microscopy.Internal organization of cells:Prokaryotic and eukaryotic cells; structure and transport through the membranes; intracelular compartments (nucleus, rough and smooth endoplasmic reticulum, golgi complex; lysosomes and peroxisomes) and protein addressing; intracellular traffic of vesicles (secretory and endocytic pathway); energy conversion (mitochondria and chloroplast); cellular communication and signalization; cytoskeleton; cellular cycle and division (mitosis and meiosis); extracellular matrix and vegetal cell wall.'
$ws.Rows.Item(18).RowHeight = 120

# Row 19
$ws.Range("A200").Copy($ws.Range("A19"))
$ws.Range("A19").Value = 'Avaliação:'

# Row 20
$ws.Range("A200").Copy($ws.Range("A20"))
$ws.Range("A20").Value = 'Método:'
$ws.Range("B200").Copy($ws.Range("B20"))
$ws.Range("B20").Value = 'Duas provas escritas (P1 e P2) distribuídas no semestre.'
$ws.Range("C200").Copy($ws.Range("C20"))
$ws.Range("C20").Value = 'Duas provas escritas (P1 e P2) distribuídas no semestre.'
$ws.Rows.Item(20).RowHeight = 60

# Row 21
$ws.Range("A200").Copy($ws.Range("A21"))
$ws.Range("A21").Value = 'Critério:'
$ws.Range("B200").Copy($ws.Range("B21"))
$ws.Range("B21").Value = 'MF=Média finalMF = (P1 + P2) / 2'
$ws.Range("C200").Copy($ws.Range("C21"))
$ws.Range("C21").Value = 'MF=Média finalMF = (P1 + P2) / 2'
$ws.Rows.Item(21).RowHeight = 60

# Row 22
$ws.Range("A200").Copy($ws.Range("A22"))
$ws.Range("A22").Value = 'Norma de recuperação:'
$ws.Range("B200").Copy($ws.Range("B22"))
$ws.Range("B22").Value = 'Nota final (NF)NF = (MF + PR)/2, onde PR é uma prova de recuperação. Prova de recuperação (PR) para alunos com Média Final maior ou igual a 3,0 e menor do que 5,0. Será considerado aprovado o aluno que tenha obtido Nota Final igual ou maior do que 5,0.'
$ws.Range("C200").Copy($ws.Range("C22"))
$ws.Range("C22").Value = 'Nota final (NF)NF = (MF + PR)/2, onde PR é uma prova de recuperação. Prova de recuperação (PR) para alunos com Média Final maior ou igual a 3,0 e menor do que 5,0. Será considerado aprovado o aluno que tenha obtido Nota Final igual ou maior do que 5,0.'
$ws.Rows.Item(22).RowHeight = 60

# Row 23
$ws.Range("A200").Copy($ws.Range("A23"))
$ws.Range("A23").Value = 'Bibliografia:'
$ws.Range("B200").Copy($ws.Range("B23"))
$ws.Range("B23").Value = '- Alberts, B., et al. Biologia Molecular da Célula, 5ed. Artmed Editora Ltda, 2010 - Cooper, G. M., Robert, E.H. A célula: uma abordagem molecular. Artmed Editora Ltda, 3° Ed. 2007.- Wasserman, S.A.; Minorsky, P.V.; Jackson, R.; Reece, J.; Cain, M.; Urry, L. Biologia de Campbell. Artmed Editora. 8 a  Edição. 2010. - Lodish, H.; Berk, A.; Matsudaira, P.; Kaiser, C. A.; Krieger, M.; Scott, M. P.; Zipurky; Darnell. Biologia Celular e Molecular. 5ª Edição. Editora Artmed, 2005.- Raven, P. H., Evert, S. E. Biologia vegetal. Editora Guanabara Koogan, 2007.- Madigan, M.T.; Martinko, J.M.; Bender, K.S.; Buckley, D.H.; Stahl, D.A. Microbiologia de Brock. Editora Artmed, 14a Edição, 2016. - Brock, T. D. ; Madigan, M.T.; Martinko, J.M.; Dunlap, P.; Clark, D. Biology of Microorganisms. Pearson Education.12a  Edição. 2009.- De Roberts, E.M.F.; Hibs, J. Bases da biologia celular e molecular. Editora Guanabara Koogan, 2006.'
$ws.Range("C200").Copy($ws.Range("C23"))
$ws.Range("C23").Value = '- Alberts, B., et al. Biologia Molecular da Célula, 5ed. Artmed Editora Ltda, 2010 - Cooper, G. M., Robert, E.H. A célula: uma abordagem molecular. Artmed Editora Ltda, 3° Ed. 2007.- Wasserman, S.A.; Minorsky, P.V.; Jackson, R.; Reece, J.; Cain, M.; Urry, L. Biologia de Campbell. Artmed Editora. 8 a  Edição. 2010. - Lodish, H.; Berk, A.; Matsudaira, P.; Kaiser, C. A.; Krieger, M.; Scott, M. P.; Zipurky; Darnell. Biologia Celular e Molecular. 5ª Edição. Editora Artmed, 2005.- Raven, P. H., Evert, S. E. Biologia vegetal. Editora Guanabara Koogan, 2007.- Madigan, M.T.; Martinko, J.M.; Bender, K.S.; Buckley, D.H.; Stahl, D.A. Microbiologia de Brock. Editora Artmed, 14a Edição, 2016. - Brock, T. D. ; Madigan, M.T.; Martinko, J.M.; Dunlap, P.; Clark, D. Biology of Microorganisms. Pearson Education.12a  Edição. 2009.- De Roberts, E.M.F.; Hibs, J. Bases da biologia celular e molecular. Editora Guanabara Koogan, 2006.'
$ws.Rows.Item(23).RowHeight = 120

# --- Clean up the scratch template row ---
$ws.Range("A200:C200").Clear()
$ws.Rows.Item(200).AutoFit()

# --- Restore selection to A1 (cosmetic, matches a freshly opened sheet) ---
$ws.Range("A1").Select()
